# Actualización automática desde tarea programada
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2's timestamp value (new run's computed value, same cell style/format retained)
$ws.Range("A2").Value = 45867.04195311342

# Append new row 3 with the latest sensor reading
$ws.Range("A3").Value = 45867.08370511012
$ws.Range("B3").Value = 2025
$ws.Range("C3").Value = 31
$ws.Range("D3").Value = 12.94
$ws.Range("E3").Value = 89.09999999999999
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = "-"
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = "02:00:32"

# Match A3's style/number format to A2 (the date/time column style)
$ws.Range("A3").NumberFormat = $ws.Range("A2").NumberFormat
